$d = $word.ActiveDocument

function Find-ParagraphContaining($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# --- 1. Title paragraph: "Lab Assignment 6" -> "Lab Assignment 5",
#        and (re)plant the "_GoBack" bookmark right after the title run,
#        mirroring Word re-stamping _GoBack at the most recent edit point. ---
$titlePara = Find-ParagraphContaining $d "Lab Assignment 6"
$titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="000B4204" w:rsidRPr="00FA22D3" w:rsidRDefault="000B4204" w:rsidP="000B4204"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:b/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr></w:pPr><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:b/><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>Lab Assignment 5</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$titlePara.Range.InsertXML($titleXml)

# --- 2. "Registration page" paragraph: collapse the five runs describing the
#        Registration page behaviour into a single run with the combined text. ---
$regPara = Find-ParagraphContaining $d "Registration page contains"
$regXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="000B4204" w:rsidRDefault="000B4204" w:rsidP="000B4204"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">The Registration page contains First name, last name, Phone number and Address fields. Based on the location, the address get populated in the fields of the registration page when the screen loads. On Clicking on the camera button camera will be opened and we can capture the photo. After capturing the photo it will be displayed in the place of camera button. On clicking on the sign in button it will navigate to the Map page.</w:t></w:r></w:p>'
$regPara.Range.InsertXML($regXml)

# --- 3. Default header: drop the stale "_GoBack" bookmark that used to mark
#        the prior edit location there. ---
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$hdrXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00773A13" w:rsidRPr="00773A13" w:rsidRDefault="00773A13" w:rsidP="00773A13"><w:pPr><w:pStyle w:val="Header"/></w:pPr></w:p>'
$hdr.Range.InsertXML($hdrXml)
